# Applies: rename all 50 sheets and update regression coefficient/p-value cells
$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ19554044",
    "summ19652918",
    "summ19754641",
    "summ19856736",
    "summ19956720",
    "summ20075842",
    "summ20203688",
    "summ20322750",
    "summ20460863",
    "summ20589505",
    "summ20705331",
    "summ20839323",
    "summ20981386",
    "summ21108317",
    "summ21240762",
    "summ21372246",
    "summ21508398",
    "summ21624401",
    "summ21757758",
    "summ21903658",
    "summ22058319",
    "summ22195273",
    "summ22325613",
    "summ22475064",
    "summ22607364",
    "summ22873207",
    "summ23005516",
    "summ23151241",
    "summ23276671",
    "summ23418156",
    "summ23541684",
    "summ23673913",
    "summ23806125",
    "summ23938914",
    "summ24070290",
    "summ24204598",
    "summ24338918",
    "summ24455463",
    "summ24590561",
    "summ24723186",
    "summ24856282",
    "summ24973003",
    "summ25108253",
    "summ25240092",
    "summ25372958",
    "summ25504604",
    "summ25622483",
    "summ25763023",
    "summ25891808",
    "summ26006618"
)

for ($i = 1; $i -le $newNames.Count; $i++) {
    $wb.Worksheets.Item($i).Name = $newNames[$i - 1]
}

# Sheet index 1
$ws = $wb.Worksheets.Item(1)
$ws.Range("B2").Value = [double]"3225.593615430962"
$ws.Range("C2").Value = [double]"1.207320548589355e-06"
$ws.Range("B3").Value = [double]"242.3401018665771"
$ws.Range("C3").Value = [double]"2.223685782038986e-08"
$ws.Range("B4").Value = [double]"-0.1010057388550846"
$ws.Range("C4").Value = [double]"0.001057814211580296"
$ws.Range("B5").Value = [double]"9678.078099375958"
$ws.Range("C5").Value = [double]"3.827646182393635e-08"
$ws.Range("B6").Value = [double]"-486.8956553057268"
$ws.Range("C6").Value = [double]"0.5480887780002639"
$ws.Range("B7").Value = [double]"-17.91152376953969"
$ws.Range("C7").Value = [double]"0.6929191661873422"

# Sheet index 2
$ws = $wb.Worksheets.Item(2)
$ws.Range("B2").Value = [double]"2745.855439874023"
$ws.Range("C2").Value = [double]"2.656633115230126e-05"
$ws.Range("B3").Value = [double]"270.6532687767415"
$ws.Range("C3").Value = [double]"2.821501414991357e-09"
$ws.Range("B4").Value = [double]"-0.06889629701065897"
$ws.Range("C4").Value = [double]"0.01880274744569366"
$ws.Range("B5").Value = [double]"10705.18800179269"
$ws.Range("C5").Value = [double]"1.266607657929519e-10"
$ws.Range("B6").Value = [double]"-784.3667059743705"
$ws.Range("C6").Value = [double]"0.3383085579784745"
$ws.Range("B7").Value = [double]"-57.86681756490349"
$ws.Range("C7").Value = [double]"0.2038414570502701"

# Sheet index 3
$ws = $wb.Worksheets.Item(3)
$ws.Range("B2").Value = [double]"4118.485915184633"
$ws.Range("C2").Value = [double]"3.10445487533501e-09"
$ws.Range("B3").Value = [double]"246.1450199267337"
$ws.Range("C3").Value = [double]"4.802683398974288e-09"
$ws.Range("B4").Value = [double]"-0.05614417245338459"
$ws.Range("C4").Value = [double]"0.03217862118894024"
$ws.Range("B5").Value = [double]"3169.380194184631"
$ws.Range("C5").Value = [double]"0.07109260014994473"
$ws.Range("B6").Value = [double]"-1303.93891123057"
$ws.Range("C6").Value = [double]"0.07377138841949302"
$ws.Range("B7").Value = [double]"14.19411394019075"
$ws.Range("C7").Value = [double]"0.7515156389461192"

# Sheet index 4
$ws = $wb.Worksheets.Item(4)
$ws.Range("B2").Value = [double]"2415.505863317058"
$ws.Range("C2").Value = [double]"5.819266274659636e-05"
$ws.Range("B3").Value = [double]"236.192909699253"
$ws.Range("C3").Value = [double]"1.05834468728579e-09"
$ws.Range("B4").Value = [double]"-0.08502472433594833"
$ws.Range("C4").Value = [double]"0.001508381309251009"
$ws.Range("B5").Value = [double]"11595.37948607568"
$ws.Range("C5").Value = [double]"8.557013055321754e-14"
$ws.Range("B6").Value = [double]"-251.8946341135891"
$ws.Range("C6").Value = [double]"0.7183677660830851"
$ws.Range("B7").Value = [double]"17.9265550988043"
$ws.Range("C7").Value = [double]"0.6558844684953573"

# Sheet index 5
$ws = $wb.Worksheets.Item(5)
$ws.Range("B2").Value = [double]"2601.420140656228"
$ws.Range("C2").Value = [double]"4.023721935002771e-05"
$ws.Range("B3").Value = [double]"244.4457769217341"
$ws.Range("C3").Value = [double]"5.165584413351316e-08"
$ws.Range("B4").Value = [double]"-0.07616546251586061"
$ws.Range("C4").Value = [double]"0.007626898551848252"
$ws.Range("B5").Value = [double]"10349.69117670941"
$ws.Range("C5").Value = [double]"3.163839585974175e-10"
$ws.Range("B6").Value = [double]"-330.7020625804109"
$ws.Range("C6").Value = [double]"0.6822936608284671"
$ws.Range("B7").Value = [double]"22.33930729746164"
$ws.Range("C7").Value = [double]"0.6313924650378087"

# Sheet index 6
$ws = $wb.Worksheets.Item(6)
$ws.Range("B2").Value = [double]"3541.911334878429"
$ws.Range("C2").Value = [double]"3.152972627070634e-09"
$ws.Range("B3").Value = [double]"245.3050506230041"
$ws.Range("C3").Value = [double]"4.8906885024039e-12"
$ws.Range("B4").Value = [double]"-0.08750248477228324"
$ws.Range("C4").Value = [double]"0.0001844554332717361"
$ws.Range("B5").Value = [double]"7212.924546857772"
$ws.Range("C5").Value = [double]"4.0094837407867e-05"
$ws.Range("B6").Value = [double]"-543.2761127334197"
$ws.Range("C6").Value = [double]"0.4116355566705352"
$ws.Range("B7").Value = [double]"-41.04136858893342"
$ws.Range("C7").Value = [double]"0.2475138680208658"

# Sheet index 7
$ws = $wb.Worksheets.Item(7)
$ws.Range("B2").Value = [double]"2360.296763405146"
$ws.Range("C2").Value = [double]"0.0002886418418439716"
$ws.Range("B3").Value = [double]"269.981579961097"
$ws.Range("C3").Value = [double]"8.194604703119572e-09"
$ws.Range("B4").Value = [double]"-0.06765081864908577"
$ws.Range("C4").Value = [double]"0.02272514710455257"
$ws.Range("B5").Value = [double]"10378.27989999433"
$ws.Range("C5").Value = [double]"2.444911685803453e-10"
$ws.Range("B6").Value = [double]"-329.5932873776376"
$ws.Range("C6").Value = [double]"0.6725053443084388"
$ws.Range("B7").Value = [double]"10.03835290804948"
$ws.Range("C7").Value = [double]"0.8348477086790516"

# Sheet index 8
$ws = $wb.Worksheets.Item(8)
$ws.Range("B2").Value = [double]"2924.199192776506"
$ws.Range("C2").Value = [double]"2.737622443926738e-05"
$ws.Range("B3").Value = [double]"244.6198196734343"
$ws.Range("C3").Value = [double]"4.714637696660897e-08"
$ws.Range("B4").Value = [double]"-0.07392562487987477"
$ws.Range("C4").Value = [double]"0.01220083971917674"
$ws.Range("B5").Value = [double]"9189.911050335766"
$ws.Range("C5").Value = [double]"3.064431109317824e-08"
$ws.Range("B6").Value = [double]"-548.7122522433635"
$ws.Range("C6").Value = [double]"0.4954804320382832"
$ws.Range("B7").Value = [double]"17.50127807982253"
$ws.Range("C7").Value = [double]"0.7112927926959092"

# Sheet index 9
$ws = $wb.Worksheets.Item(9)
$ws.Range("B2").Value = [double]"3045.057527541785"
$ws.Range("C2").Value = [double]"7.427005969574003e-06"
$ws.Range("B3").Value = [double]"246.1738882181833"
$ws.Range("C3").Value = [double]"4.133226052058102e-08"
$ws.Range("B4").Value = [double]"-0.07237828153025164"
$ws.Range("C4").Value = [double]"0.01714959521179016"
$ws.Range("B5").Value = [double]"10212.93833415287"
$ws.Range("C5").Value = [double]"1.519828055652078e-09"
$ws.Range("B6").Value = [double]"-1126.604140811113"
$ws.Range("C6").Value = [double]"0.1865797015347208"
$ws.Range("B7").Value = [double]"-4.028425812093801"
$ws.Range("C7").Value = [double]"0.9327583832242931"

# Sheet index 10
$ws = $wb.Worksheets.Item(10)
$ws.Range("B2").Value = [double]"3126.311552229487"
$ws.Range("C2").Value = [double]"2.141159238289316e-06"
$ws.Range("B3").Value = [double]"240.213525469004"
$ws.Range("C3").Value = [double]"6.796354843859411e-08"
$ws.Range("B4").Value = [double]"-0.08616378604976196"
$ws.Range("C4").Value = [double]"0.004102729209966388"
$ws.Range("B5").Value = [double]"9471.260654983891"
$ws.Range("C5").Value = [double]"3.464605924737363e-08"
$ws.Range("B6").Value = [double]"-659.382469008264"
$ws.Range("C6").Value = [double]"0.4107317685729223"
$ws.Range("B7").Value = [double]"-13.75929446180182"
$ws.Range("C7").Value = [double]"0.7647840549022381"

# Sheet index 11
$ws = $wb.Worksheets.Item(11)
$ws.Range("B2").Value = [double]"2669.822636102742"
$ws.Range("C2").Value = [double]"7.06138434644017e-05"
$ws.Range("B3").Value = [double]"270.3726102123831"
$ws.Range("C3").Value = [double]"6.507567483856899e-09"
$ws.Range("B4").Value = [double]"-0.06656862753282161"
$ws.Range("C4").Value = [double]"0.02711563352911538"
$ws.Range("B5").Value = [double]"10217.43869129972"
$ws.Range("C5").Value = [double]"2.862648855800052e-09"
$ws.Range("B6").Value = [double]"-688.1475695989354"
$ws.Range("C6").Value = [double]"0.389335365668203"
$ws.Range("B7").Value = [double]"-15.60150837557268"
$ws.Range("C7").Value = [double]"0.7407932956986019"

# Sheet index 12
$ws = $wb.Worksheets.Item(12)
$ws.Range("B2").Value = [double]"2719.99547371413"
$ws.Range("C2").Value = [double]"5.248644658618557e-05"
$ws.Range("B3").Value = [double]"267.3828037949182"
$ws.Range("C3").Value = [double]"7.701084223326059e-09"
$ws.Range("B4").Value = [double]"-0.06595849826948427"
$ws.Range("C4").Value = [double]"0.03485595782077895"
$ws.Range("B5").Value = [double]"9179.50379110569"
$ws.Range("C5").Value = [double]"2.711850574120844e-08"
$ws.Range("B6").Value = [double]"-744.0729301262165"
$ws.Range("C6").Value = [double]"0.3792184234749619"
$ws.Range("B7").Value = [double]"34.51756003992065"
$ws.Range("C7").Value = [double]"0.4796949142892662"

# Sheet index 13
$ws = $wb.Worksheets.Item(13)
$ws.Range("B2").Value = [double]"3121.881020528695"
$ws.Range("C2").Value = [double]"7.103600953621502e-06"
$ws.Range("B3").Value = [double]"244.0944499910362"
$ws.Range("C3").Value = [double]"9.849139323830115e-08"
$ws.Range("B4").Value = [double]"-0.08970083891304134"
$ws.Range("C4").Value = [double]"0.003676461169279479"
$ws.Range("B5").Value = [double]"9751.000097418411"
$ws.Range("C5").Value = [double]"7.861612598127544e-09"
$ws.Range("B6").Value = [double]"-559.1750817613222"
$ws.Range("C6").Value = [double]"0.501081553455507"
$ws.Range("B7").Value = [double]"-24.80950474991621"
$ws.Range("C7").Value = [double]"0.5988274756465208"

# Sheet index 14
$ws = $wb.Worksheets.Item(14)
$ws.Range("B2").Value = [double]"4019.083358505906"
$ws.Range("C2").Value = [double]"3.720922567093477e-09"
$ws.Range("B3").Value = [double]"231.8044430607773"
$ws.Range("C3").Value = [double]"1.733013593983474e-08"
$ws.Range("B4").Value = [double]"-0.08734530285064404"
$ws.Range("C4").Value = [double]"0.001197036290905893"
$ws.Range("B5").Value = [double]"5819.159605652239"
$ws.Range("C5").Value = [double]"0.002125934715468424"
$ws.Range("B6").Value = [double]"-829.3809700213301"
$ws.Range("C6").Value = [double]"0.2596100000589794"
$ws.Range("B7").Value = [double]"-43.31850142843673"
$ws.Range("C7").Value = [double]"0.2937865228428678"

# Sheet index 15
$ws = $wb.Worksheets.Item(15)
$ws.Range("B2").Value = [double]"2484.536526379713"
$ws.Range("C2").Value = [double]"5.339340504675882e-06"
$ws.Range("B3").Value = [double]"235.6684323833725"
$ws.Range("C3").Value = [double]"4.633522057927032e-11"
$ws.Range("B4").Value = [double]"-0.07612887217435785"
$ws.Range("C4").Value = [double]"0.001102579608012477"
$ws.Range("B5").Value = [double]"11317.18101075725"
$ws.Range("C5").Value = [double]"3.060215547527741e-14"
$ws.Range("B6").Value = [double]"-454.2380947982674"
$ws.Range("C6").Value = [double]"0.498218312533949"
$ws.Range("B7").Value = [double]"15.29738684699623"
$ws.Range("C7").Value = [double]"0.6929937717956505"

# Sheet index 16
$ws = $wb.Worksheets.Item(16)
$ws.Range("B2").Value = [double]"3008.645653747376"
$ws.Range("C2").Value = [double]"1.615967911502888e-05"
$ws.Range("B3").Value = [double]"244.4793320047651"
$ws.Range("C3").Value = [double]"1.334441948022671e-07"
$ws.Range("B4").Value = [double]"-0.07689737307303875"
$ws.Range("C4").Value = [double]"0.01532211808875623"
$ws.Range("B5").Value = [double]"9390.411306405291"
$ws.Range("C5").Value = [double]"4.91216869555938e-08"
$ws.Range("B6").Value = [double]"-643.5858732228246"
$ws.Range("C6").Value = [double]"0.4433735699651856"
$ws.Range("B7").Value = [double]"-0.2990292690563905"
$ws.Range("C7").Value = [double]"0.9949690686851964"

# Sheet index 17
$ws = $wb.Worksheets.Item(17)
$ws.Range("B2").Value = [double]"2989.192707166097"
$ws.Range("C2").Value = [double]"6.922166659216829e-06"
$ws.Range("B3").Value = [double]"248.5512491880194"
$ws.Range("C3").Value = [double]"6.55760584547614e-08"
$ws.Range("B4").Value = [double]"-0.08193179620630706"
$ws.Range("C4").Value = [double]"0.006269623753543351"
$ws.Range("B5").Value = [double]"10415.93709848066"
$ws.Range("C5").Value = [double]"2.98597821037694e-09"
$ws.Range("B6").Value = [double]"-732.420901419839"
$ws.Range("C6").Value = [double]"0.379985021642909"
$ws.Range("B7").Value = [double]"-14.47534933874883"
$ws.Range("C7").Value = [double]"0.7647946454822598"

# Sheet index 18
$ws = $wb.Worksheets.Item(18)
$ws.Range("B2").Value = [double]"3296.522210554763"
$ws.Range("C2").Value = [double]"3.940525235656022e-09"
$ws.Range("B3").Value = [double]"256.8079768820498"
$ws.Range("C3").Value = [double]"1.799459501629301e-13"
$ws.Range("B4").Value = [double]"-0.07989250536157506"
$ws.Range("C4").Value = [double]"0.0005447355861443601"
$ws.Range("B5").Value = [double]"7001.044921844376"
$ws.Range("C5").Value = [double]"4.97984293940687e-06"
$ws.Range("B6").Value = [double]"-400.8053470290711"
$ws.Range("C6").Value = [double]"0.5058884776891683"
$ws.Range("B7").Value = [double]"-53.45060744615072"
$ws.Range("C7").Value = [double]"0.1162018606399586"

# Sheet index 19
$ws = $wb.Worksheets.Item(19)
$ws.Range("B2").Value = [double]"3088.922423645701"
$ws.Range("C2").Value = [double]"4.881909068467552e-06"
$ws.Range("B3").Value = [double]"238.0398614084633"
$ws.Range("C3").Value = [double]"1.267217573662908e-07"
$ws.Range("B4").Value = [double]"-0.06822379830167527"
$ws.Range("C4").Value = [double]"0.02096780922163937"
$ws.Range("B5").Value = [double]"8859.000922105331"
$ws.Range("C5").Value = [double]"1.415565789083799e-07"
$ws.Range("B6").Value = [double]"-1037.703321452315"
$ws.Range("C6").Value = [double]"0.2040952447593855"
$ws.Range("B7").Value = [double]"38.71105224236697"
$ws.Range("C7").Value = [double]"0.4160192680940968"

# Sheet index 20
$ws = $wb.Worksheets.Item(20)
$ws.Range("B2").Value = [double]"2652.204537520405"
$ws.Range("C2").Value = [double]"5.237372590348826e-05"
$ws.Range("B3").Value = [double]"252.1118012056651"
$ws.Range("C3").Value = [double]"6.33237428118254e-09"
$ws.Range("B4").Value = [double]"-0.0831501712569027"
$ws.Range("C4").Value = [double]"0.003054654696692422"
$ws.Range("B5").Value = [double]"10954.79452924"
$ws.Range("C5").Value = [double]"3.156902910880581e-11"
$ws.Range("B6").Value = [double]"-415.0745534639277"
$ws.Range("C6").Value = [double]"0.601481343743739"
$ws.Range("B7").Value = [double]"2.264041008344932"
$ws.Range("C7").Value = [double]"0.9612271064767628"

# Sheet index 21
$ws = $wb.Worksheets.Item(21)
$ws.Range("B2").Value = [double]"2700.262952625396"
$ws.Range("C2").Value = [double]"2.65377878535797e-05"
$ws.Range("B3").Value = [double]"243.305233645691"
$ws.Range("C3").Value = [double]"4.582061786008416e-08"
$ws.Range("B4").Value = [double]"-0.07857088535214279"
$ws.Range("C4").Value = [double]"0.004940020741907983"
$ws.Range("B5").Value = [double]"10008.51833127044"
$ws.Range("C5").Value = [double]"2.034986919187708e-09"
$ws.Range("B6").Value = [double]"-390.9661218437523"
$ws.Range("C6").Value = [double]"0.6246904223290978"
$ws.Range("B7").Value = [double]"20.95329082175201"
$ws.Range("C7").Value = [double]"0.6636108181218525"

# Sheet index 22
$ws = $wb.Worksheets.Item(22)
$ws.Range("B2").Value = [double]"2812.466621587015"
$ws.Range("C2").Value = [double]"2.369429630937393e-05"
$ws.Range("B3").Value = [double]"242.6337800971718"
$ws.Range("C3").Value = [double]"1.777670920527133e-08"
$ws.Range("B4").Value = [double]"-0.06816128327477509"
$ws.Range("C4").Value = [double]"0.01930204934155853"
$ws.Range("B5").Value = [double]"10153.05808415007"
$ws.Range("C5").Value = [double]"1.207174352808371e-09"
$ws.Range("B6").Value = [double]"-810.9544267061705"
$ws.Range("C6").Value = [double]"0.3084491563438367"
$ws.Range("B7").Value = [double]"9.901670107799831"
$ws.Range("C7").Value = [double]"0.8274351728044589"

# Sheet index 23
$ws = $wb.Worksheets.Item(23)
$ws.Range("B2").Value = [double]"4269.051016014859"
$ws.Range("C2").Value = [double]"7.534785391132907e-12"
$ws.Range("B3").Value = [double]"224.9486237877997"
$ws.Range("C3").Value = [double]"4.732137575460817e-10"
$ws.Range("B4").Value = [double]"-0.07484446544020863"
$ws.Range("C4").Value = [double]"0.001648096997863937"
$ws.Range("B5").Value = [double]"4381.945800337097"
$ws.Range("C5").Value = [double]"0.008341261425867065"
$ws.Range("B6").Value = [double]"-1146.670089229862"
$ws.Range("C6").Value = [double]"0.07677099707180687"
$ws.Range("B7").Value = [double]"-34.12124351175093"
$ws.Range("C7").Value = [double]"0.3413543847509274"

# Sheet index 24
$ws = $wb.Worksheets.Item(24)
$ws.Range("B2").Value = [double]"2637.641396573023"
$ws.Range("C2").Value = [double]"8.045261603626105e-05"
$ws.Range("B3").Value = [double]"269.8717393647566"
$ws.Range("C3").Value = [double]"3.166417576584487e-09"
$ws.Range("B4").Value = [double]"-0.07473582243966388"
$ws.Range("C4").Value = [double]"0.02129813334813782"
$ws.Range("B5").Value = [double]"9832.968816331795"
$ws.Range("C5").Value = [double]"2.022072819292146e-09"
$ws.Range("B6").Value = [double]"-443.4360197610072"
$ws.Range("C6").Value = [double]"0.5821038276750052"
$ws.Range("B7").Value = [double]"11.31990793943604"
$ws.Range("C7").Value = [double]"0.8120428492607443"

# Sheet index 25
$ws = $wb.Worksheets.Item(25)
$ws.Range("B2").Value = [double]"2918.387645489185"
$ws.Range("C2").Value = [double]"1.651192115413219e-05"
$ws.Range("B3").Value = [double]"257.4067264883326"
$ws.Range("C3").Value = [double]"2.927820091687992e-08"
$ws.Range("B4").Value = [double]"-0.08945711138911197"
$ws.Range("C4").Value = [double]"0.002261451225937799"
$ws.Range("B5").Value = [double]"11077.37515232904"
$ws.Range("C5").Value = [double]"2.392427737467984e-10"
$ws.Range("B6").Value = [double]"-624.6010281179204"
$ws.Range("C6").Value = [double]"0.456749120898199"
$ws.Range("B7").Value = [double]"-43.65966950426267"
$ws.Range("C7").Value = [double]"0.3459537046890921"

# Sheet index 26
$ws = $wb.Worksheets.Item(26)
$ws.Range("B2").Value = [double]"4351.524372924243"
$ws.Range("C2").Value = [double]"6.227777949621776e-11"
$ws.Range("B3").Value = [double]"213.110524228927"
$ws.Range("C3").Value = [double]"1.354902027913702e-07"
$ws.Range("B4").Value = [double]"-0.08452536525358895"
$ws.Range("C4").Value = [double]"0.001405146353849094"
$ws.Range("B5").Value = [double]"4111.099861331712"
$ws.Range("C5").Value = [double]"0.02270153993976514"
$ws.Range("B6").Value = [double]"-961.301748357529"
$ws.Range("C6").Value = [double]"0.1870759801500239"
$ws.Range("B7").Value = [double]"-23.25465259661172"
$ws.Range("C7").Value = [double]"0.588247648629187"

# Sheet index 27
$ws = $wb.Worksheets.Item(27)
$ws.Range("B2").Value = [double]"2775.853016183327"
$ws.Range("C2").Value = [double]"2.878870821265456e-05"
$ws.Range("B3").Value = [double]"261.5250688184043"
$ws.Range("C3").Value = [double]"8.876772646014763e-09"
$ws.Range("B4").Value = [double]"-0.06631710934582521"
$ws.Range("C4").Value = [double]"0.02166562556258072"
$ws.Range("B5").Value = [double]"10570.36438049148"
$ws.Range("C5").Value = [double]"1.377411337745951e-10"
$ws.Range("B6").Value = [double]"-984.4053600534007"
$ws.Range("C6").Value = [double]"0.2295258543103504"
$ws.Range("B7").Value = [double]"-16.52080133237216"
$ws.Range("C7").Value = [double]"0.7258527204946648"

# Sheet index 28
$ws = $wb.Worksheets.Item(28)
$ws.Range("B2").Value = [double]"3430.230358656257"
$ws.Range("C2").Value = [double]"9.493720871191482e-07"
$ws.Range("B3").Value = [double]"228.5179837358748"
$ws.Range("C3").Value = [double]"3.68688407441291e-07"
$ws.Range("B4").Value = [double]"-0.08700186385584677"
$ws.Range("C4").Value = [double]"0.004038991602751642"
$ws.Range("B5").Value = [double]"8915.227372384454"
$ws.Range("C5").Value = [double]"7.755213549589143e-07"
$ws.Range("B6").Value = [double]"-851.1697664699329"
$ws.Range("C6").Value = [double]"0.3009598766036882"
$ws.Range("B7").Value = [double]"7.585865449817618"
$ws.Range("C7").Value = [double]"0.8707121162473392"

# Sheet index 29
$ws = $wb.Worksheets.Item(29)
$ws.Range("B2").Value = [double]"1855.208451801971"
$ws.Range("C2").Value = [double]"0.0005922332621460852"
$ws.Range("B3").Value = [double]"273.7455925848309"
$ws.Range("C3").Value = [double]"1.261575718884033e-12"
$ws.Range("B4").Value = [double]"-0.08059228102886891"
$ws.Range("C4").Value = [double]"0.001184696320197894"
$ws.Range("B5").Value = [double]"12169.25009080665"
$ws.Range("C5").Value = [double]"4.737269223106379e-17"
$ws.Range("B6").Value = [double]"312.4400013488359"
$ws.Range("C6").Value = [double]"0.6343676330247205"
$ws.Range("B7").Value = [double]"-16.47062344409717"
$ws.Range("C7").Value = [double]"0.6566659332706442"

# Sheet index 30
$ws = $wb.Worksheets.Item(30)
$ws.Range("B2").Value = [double]"2777.018269946776"
$ws.Range("C2").Value = [double]"4.285123057287583e-05"
$ws.Range("B3").Value = [double]"264.8422664635159"
$ws.Range("C3").Value = [double]"2.479486032343143e-09"
$ws.Range("B4").Value = [double]"-0.07099847548889299"
$ws.Range("C4").Value = [double]"0.02128788525016232"
$ws.Range("B5").Value = [double]"9883.611110656031"
$ws.Range("C5").Value = [double]"4.451903777583983e-09"
$ws.Range("B6").Value = [double]"-732.4322973952242"
$ws.Range("C6").Value = [double]"0.377413977794671"
$ws.Range("B7").Value = [double]"5.080494347825066"
$ws.Range("C7").Value = [double]"0.9152070130271373"

# Sheet index 31
$ws = $wb.Worksheets.Item(31)
$ws.Range("B2").Value = [double]"2577.094917741621"
$ws.Range("C2").Value = [double]"4.277293291377556e-05"
$ws.Range("B3").Value = [double]"243.2161373196716"
$ws.Range("C3").Value = [double]"1.586698075843286e-08"
$ws.Range("B4").Value = [double]"-0.07790622534311459"
$ws.Range("C4").Value = [double]"0.00529219343419115"
$ws.Range("B5").Value = [double]"11091.87904739257"
$ws.Range("C5").Value = [double]"2.425914983185929e-11"
$ws.Range("B6").Value = [double]"-402.9890751312018"
$ws.Range("C6").Value = [double]"0.6053539536348704"
$ws.Range("B7").Value = [double]"-17.9507293013587"
$ws.Range("C7").Value = [double]"0.6897078999217101"

# Sheet index 32
$ws = $wb.Worksheets.Item(32)
$ws.Range("B2").Value = [double]"3057.882084446982"
$ws.Range("C2").Value = [double]"4.223768364722519e-06"
$ws.Range("B3").Value = [double]"238.5053456430184"
$ws.Range("C3").Value = [double]"1.889009961354052e-07"
$ws.Range("B4").Value = [double]"-0.08926880357542472"
$ws.Range("C4").Value = [double]"0.004229220780051073"
$ws.Range("B5").Value = [double]"9917.739485678712"
$ws.Range("C5").Value = [double]"6.509185485610062e-09"
$ws.Range("B6").Value = [double]"-776.5691319396578"
$ws.Range("C6").Value = [double]"0.348418122555181"
$ws.Range("B7").Value = [double]"-6.676131620496065"
$ws.Range("C7").Value = [double]"0.8881780241068967"

# Sheet index 33
$ws = $wb.Worksheets.Item(33)
$ws.Range("B2").Value = [double]"3573.677172692706"
$ws.Range("C2").Value = [double]"2.256847285302772e-09"
$ws.Range("B3").Value = [double]"259.6275578617732"
$ws.Range("C3").Value = [double]"9.90564024518357e-13"
$ws.Range("B4").Value = [double]"-0.07607495599482772"
$ws.Range("C4").Value = [double]"0.001044239368421336"
$ws.Range("B5").Value = [double]"5470.901309341269"
$ws.Range("C5").Value = [double]"0.001063299889398755"
$ws.Range("B6").Value = [double]"-249.9613266185302"
$ws.Range("C6").Value = [double]"0.6925684169244373"
$ws.Range("B7").Value = [double]"-49.91141827653721"
$ws.Range("C7").Value = [double]"0.1857578474325319"

# Sheet index 34
$ws = $wb.Worksheets.Item(34)
$ws.Range("B2").Value = [double]"2461.499667124758"
$ws.Range("C2").Value = [double]"0.0003800217921094492"
$ws.Range("B3").Value = [double]"281.018730329074"
$ws.Range("C3").Value = [double]"1.286978402418189e-09"
$ws.Range("B4").Value = [double]"-0.05674013307158164"
$ws.Range("C4").Value = [double]"0.0588770648734962"
$ws.Range("B5").Value = [double]"9873.918781615122"
$ws.Range("C5").Value = [double]"3.820768373159416e-09"
$ws.Range("B6").Value = [double]"-840.2556384750042"
$ws.Range("C6").Value = [double]"0.2912444898227867"
$ws.Range("B7").Value = [double]"49.39148714871479"
$ws.Range("C7").Value = [double]"0.3173879086368263"

# Sheet index 35
$ws = $wb.Worksheets.Item(35)
$ws.Range("B2").Value = [double]"3432.274181889614"
$ws.Range("C2").Value = [double]"5.59774499516534e-07"
$ws.Range("B3").Value = [double]"226.5381326647409"
$ws.Range("C3").Value = [double]"3.460007310718317e-07"
$ws.Range("B4").Value = [double]"-0.08416047150923285"
$ws.Range("C4").Value = [double]"0.005481355648571913"
$ws.Range("B5").Value = [double]"9745.763615152668"
$ws.Range("C5").Value = [double]"5.56202422452575e-09"
$ws.Range("B6").Value = [double]"-1146.73696451976"
$ws.Range("C6").Value = [double]"0.1755707552574358"
$ws.Range("B7").Value = [double]"-8.106560446059291"
$ws.Range("C7").Value = [double]"0.8560332981951131"

# Sheet index 36
$ws = $wb.Worksheets.Item(36)
$ws.Range("B2").Value = [double]"2842.317912592272"
$ws.Range("C2").Value = [double]"1.3161525005869e-05"
$ws.Range("B3").Value = [double]"234.5709485433441"
$ws.Range("C3").Value = [double]"5.454676998305403e-08"
$ws.Range("B4").Value = [double]"-0.08356512348758548"
$ws.Range("C4").Value = [double]"0.004761916837537638"
$ws.Range("B5").Value = [double]"11154.15362733887"
$ws.Range("C5").Value = [double]"2.594287266514939e-10"
$ws.Range("B6").Value = [double]"-633.9216271773639"
$ws.Range("C6").Value = [double]"0.4303061129358586"
$ws.Range("B7").Value = [double]"-9.952847814307829"
$ws.Range("C7").Value = [double]"0.8262378379358228"

# Sheet index 37
$ws = $wb.Worksheets.Item(37)
$ws.Range("B2").Value = [double]"4301.068544896866"
$ws.Range("C2").Value = [double]"8.518727043894071e-10"
$ws.Range("B3").Value = [double]"241.081308305941"
$ws.Range("C3").Value = [double]"2.389838196529725e-08"
$ws.Range("B4").Value = [double]"-0.09111837170799343"
$ws.Range("C4").Value = [double]"0.001321701485053344"
$ws.Range("B5").Value = [double]"4433.227399405382"
$ws.Range("C5").Value = [double]"0.01225569980715509"
$ws.Range("B6").Value = [double]"-959.968247633841"
$ws.Range("C6").Value = [double]"0.2266233120736972"
$ws.Range("B7").Value = [double]"-44.97953023982478"
$ws.Range("C7").Value = [double]"0.2907332484494225"

# Sheet index 38
$ws = $wb.Worksheets.Item(38)
$ws.Range("B2").Value = [double]"2673.475701089908"
$ws.Range("C2").Value = [double]"3.767667037607327e-05"
$ws.Range("B3").Value = [double]"245.1626261735965"
$ws.Range("C3").Value = [double]"1.799211970307246e-08"
$ws.Range("B4").Value = [double]"-0.07836259649430891"
$ws.Range("C4").Value = [double]"0.007105879705965745"
$ws.Range("B5").Value = [double]"9868.233156992919"
$ws.Range("C5").Value = [double]"4.066295135776989e-09"
$ws.Range("B6").Value = [double]"-263.2296903171157"
$ws.Range("C6").Value = [double]"0.7443333470174174"
$ws.Range("B7").Value = [double]"15.31367381246797"
$ws.Range("C7").Value = [double]"0.7414870411657709"

# Sheet index 39
$ws = $wb.Worksheets.Item(39)
$ws.Range("B2").Value = [double]"2240.477176708732"
$ws.Range("C2").Value = [double]"0.0001523206008793173"
$ws.Range("B3").Value = [double]"270.1776569868948"
$ws.Range("C3").Value = [double]"4.046741075538347e-11"
$ws.Range("B4").Value = [double]"-0.06306246382650955"
$ws.Range("C4").Value = [double]"0.01606181399891539"
$ws.Range("B5").Value = [double]"10489.26364002045"
$ws.Range("C5").Value = [double]"5.197096136852652e-12"
$ws.Range("B6").Value = [double]"-338.8380985539391"
$ws.Range("C6").Value = [double]"0.6206490447095004"
$ws.Range("B7").Value = [double]"21.15935781433398"
$ws.Range("C7").Value = [double]"0.6313615163276931"

# Sheet index 40
$ws = $wb.Worksheets.Item(40)
$ws.Range("B2").Value = [double]"3051.694321851863"
$ws.Range("C2").Value = [double]"5.434914708847623e-06"
$ws.Range("B3").Value = [double]"251.2493136835678"
$ws.Range("C3").Value = [double]"2.488548798562708e-08"
$ws.Range("B4").Value = [double]"-0.07652383353917527"
$ws.Range("C4").Value = [double]"0.008633912617925149"
$ws.Range("B5").Value = [double]"10224.09161719736"
$ws.Range("C5").Value = [double]"8.236361939815721e-10"
$ws.Range("B6").Value = [double]"-987.5741884432981"
$ws.Range("C6").Value = [double]"0.2170457986477776"
$ws.Range("B7").Value = [double]"-15.61365429666606"
$ws.Range("C7").Value = [double]"0.7335478343474573"

# Sheet index 41
$ws = $wb.Worksheets.Item(41)
$ws.Range("B2").Value = [double]"2469.541046041305"
$ws.Range("C2").Value = [double]"6.799884007203256e-05"
$ws.Range("B3").Value = [double]"253.9384351849875"
$ws.Range("C3").Value = [double]"7.257661405847422e-10"
$ws.Range("B4").Value = [double]"-0.09878407055711888"
$ws.Range("C4").Value = [double]"0.0008876652402871117"
$ws.Range("B5").Value = [double]"11021.88083621026"
$ws.Range("C5").Value = [double]"3.16884112356434e-12"
$ws.Range("B6").Value = [double]"-27.28508241568966"
$ws.Range("C6").Value = [double]"0.9703929590424376"
$ws.Range("B7").Value = [double]"-17.69387475571234"
$ws.Range("C7").Value = [double]"0.6636307147189435"

# Sheet index 42
$ws = $wb.Worksheets.Item(42)
$ws.Range("B2").Value = [double]"2492.349132188724"
$ws.Range("C2").Value = [double]"0.0002436843656911054"
$ws.Range("B3").Value = [double]"260.908212514624"
$ws.Range("C3").Value = [double]"1.069622827444285e-08"
$ws.Range("B4").Value = [double]"-0.05624895563721107"
$ws.Range("C4").Value = [double]"0.05579303037120796"
$ws.Range("B5").Value = [double]"9843.062478037276"
$ws.Range("C5").Value = [double]"1.668312936971856e-09"
$ws.Range("B6").Value = [double]"-691.9198921803068"
$ws.Range("C6").Value = [double]"0.3795922341091816"
$ws.Range("B7").Value = [double]"46.69933169531209"
$ws.Range("C7").Value = [double]"0.3305234528885812"

# Sheet index 43
$ws = $wb.Worksheets.Item(43)
$ws.Range("B2").Value = [double]"2419.000425054866"
$ws.Range("C2").Value = [double]"9.379817187419878e-05"
$ws.Range("B3").Value = [double]"258.6585984580032"
$ws.Range("C3").Value = [double]"9.829618163309115e-10"
$ws.Range("B4").Value = [double]"-0.07178596154392625"
$ws.Range("C4").Value = [double]"0.008970244758986931"
$ws.Range("B5").Value = [double]"11632.52143931545"
$ws.Range("C5").Value = [double]"1.289889903828458e-11"
$ws.Range("B6").Value = [double]"-501.0241488765005"
$ws.Range("C6").Value = [double]"0.5083479366634798"
$ws.Range("B7").Value = [double]"-14.81810526930089"
$ws.Range("C7").Value = [double]"0.7488284180083544"

# Sheet index 44
$ws = $wb.Worksheets.Item(44)
$ws.Range("B2").Value = [double]"4327.240953373301"
$ws.Range("C2").Value = [double]"1.141824319388805e-10"
$ws.Range("B3").Value = [double]"229.9780918920361"
$ws.Range("C3").Value = [double]"1.824235485132822e-08"
$ws.Range("B4").Value = [double]"-0.0853886109987541"
$ws.Range("C4").Value = [double]"0.001211615795639189"
$ws.Range("B5").Value = [double]"4284.316229769214"
$ws.Range("C5").Value = [double]"0.01528721653919161"
$ws.Range("B6").Value = [double]"-969.37903604805"
$ws.Range("C6").Value = [double]"0.2041827179012777"
$ws.Range("B7").Value = [double]"-18.43222736089021"
$ws.Range("C7").Value = [double]"0.6665977064331342"

# Sheet index 45
$ws = $wb.Worksheets.Item(45)
$ws.Range("B2").Value = [double]"3285.793981769412"
$ws.Range("C2").Value = [double]"1.001947139137276e-06"
$ws.Range("B3").Value = [double]"240.9622371131585"
$ws.Range("C3").Value = [double]"1.015436688128774e-07"
$ws.Range("B4").Value = [double]"-0.08177478300728204"
$ws.Range("C4").Value = [double]"0.006072173605681318"
$ws.Range("B5").Value = [double]"9561.156101101631"
$ws.Range("C5").Value = [double]"1.951524650664004e-08"
$ws.Range("B6").Value = [double]"-960.7186203480849"
$ws.Range("C6").Value = [double]"0.249461994802166"
$ws.Range("B7").Value = [double]"-21.57083002428234"
$ws.Range("C7").Value = [double]"0.6359371510065737"

# Sheet index 46
$ws = $wb.Worksheets.Item(46)
$ws.Range("B2").Value = [double]"3092.088569336216"
$ws.Range("C2").Value = [double]"4.371966398370198e-06"
$ws.Range("B3").Value = [double]"235.7843012365968"
$ws.Range("C3").Value = [double]"1.420969913519254e-07"
$ws.Range("B4").Value = [double]"-0.07301978321873533"
$ws.Range("C4").Value = [double]"0.02026924801511395"
$ws.Range("B5").Value = [double]"9354.332475258692"
$ws.Range("C5").Value = [double]"1.969397165131967e-08"
$ws.Range("B6").Value = [double]"-989.9559865969263"
$ws.Range("C6").Value = [double]"0.23979750281013"
$ws.Range("B7").Value = [double]"24.5862342219942"
$ws.Range("C7").Value = [double]"0.6015439489352838"

# Sheet index 47
$ws = $wb.Worksheets.Item(47)
$ws.Range("B2").Value = [double]"2419.220702920852"
$ws.Range("C2").Value = [double]"0.0002537585228407168"
$ws.Range("B3").Value = [double]"297.0451766380121"
$ws.Range("C3").Value = [double]"7.922607492591836e-11"
$ws.Range("B4").Value = [double]"-0.06262535885650711"
$ws.Range("C4").Value = [double]"0.02778373100367122"
$ws.Range("B5").Value = [double]"11160.2682324121"
$ws.Range("C5").Value = [double]"3.111425983305699e-11"
$ws.Range("B6").Value = [double]"-853.9321143963888"
$ws.Range("C6").Value = [double]"0.2856833717373911"
$ws.Range("B7").Value = [double]"-13.36407907704722"
$ws.Range("C7").Value = [double]"0.7849387165101767"

# Sheet index 48
$ws = $wb.Worksheets.Item(48)
$ws.Range("B2").Value = [double]"4056.331991337256"
$ws.Range("C2").Value = [double]"7.201405879625956e-10"
$ws.Range("B3").Value = [double]"232.7975746951775"
$ws.Range("C3").Value = [double]"2.49283839811647e-09"
$ws.Range("B4").Value = [double]"-0.08488644427584856"
$ws.Range("C4").Value = [double]"0.001588951148620734"
$ws.Range("B5").Value = [double]"4376.953192469316"
$ws.Range("C5").Value = [double]"0.01493419225583836"
$ws.Range("B6").Value = [double]"-619.3945831219636"
$ws.Range("C6").Value = [double]"0.3964372176930054"
$ws.Range("B7").Value = [double]"-27.94391927878804"
$ws.Range("C7").Value = [double]"0.5038092210447145"

# Sheet index 49
$ws = $wb.Worksheets.Item(49)
$ws.Range("B2").Value = [double]"2126.929947786504"
$ws.Range("C2").Value = [double]"0.0001768107030768924"
$ws.Range("B3").Value = [double]"262.8891690378232"
$ws.Range("C3").Value = [double]"5.249718641004489e-11"
$ws.Range("B4").Value = [double]"-0.07805001298041264"
$ws.Range("C4").Value = [double]"0.001344992912147491"
$ws.Range("B5").Value = [double]"11538.4805354267"
$ws.Range("C5").Value = [double]"2.266182912541231e-14"
$ws.Range("B6").Value = [double]"42.7128622665939"
$ws.Range("C6").Value = [double]"0.9497530012230658"
$ws.Range("B7").Value = [double]"-25.60751112891613"
$ws.Range("C7").Value = [double]"0.5163125346841999"

# Sheet index 50
$ws = $wb.Worksheets.Item(50)
$ws.Range("B2").Value = [double]"3452.645920527577"
$ws.Range("C2").Value = [double]"6.587467665932092e-07"
$ws.Range("B3").Value = [double]"212.1190766445928"
$ws.Range("C3").Value = [double]"2.786872197180397e-06"
$ws.Range("B4").Value = [double]"-0.09921946021827617"
$ws.Range("C4").Value = [double]"0.001355202188633512"
$ws.Range("B5").Value = [double]"9608.407379761684"
$ws.Range("C5").Value = [double]"3.72844457416929e-08"
$ws.Range("B6").Value = [double]"-658.0034159554834"
$ws.Range("C6").Value = [double]"0.4235832503314012"
$ws.Range("B7").Value = [double]"3.041704204087594"
$ws.Range("C7").Value = [double]"0.9464352748786223"
